$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert an extra space between first and last name in column A (rows 2-204)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $name = $cell.Value2
    if ($name -ne $null -and $name -match ' ') {
        $cell.Value = ($name -replace ' ', '  ')
    }
}

# 2. Apply AutoFilter over the data range
$dataRange = $ws.Range("A1:C204")
$dataRange.AutoFilter()

# 3. Make sure the hidden built-in _FilterDatabase defined name exists (sheet-scoped)
$fd = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$C`$204")
$fd.Visible = $false

# 4. Move the active selection to A2
$ws.Range("A2").Select()
